$wb = $excel.ActiveWorkbook

# Map of worksheet name -> hashtable of row -> new F-column value
$changes = @{
    "展览" = @{
        3  = 21323
        8  = 7952
        12 = 316
        13 = 66
        20 = 539
        27 = 1191
        30 = 227
        33 = 8
        46 = 440
    }
    "全部类型" = @{
        3  = 21323
        7  = 7952
        11 = 316
        12 = 66
        18 = 539
        25 = 1191
        28 = 227
        32 = 8
        46 = 440
    }
}

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsMap = $changes[$sheetName]
    foreach ($row in $rowsMap.Keys) {
        $value = $rowsMap[$row]
        $ws.Cells.Item([int]$row, 6).Value = $value
    }
}
